$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 9-18: the underlying data export was refreshed, causing records
# to be reshuffled among these rows, coordinates (Q/R) to be rounded to whole
# metres, and the Starttid/Sluttid time columns (Z, AB) to be cleared out.

# Row 9
$ws.Range("A9").Value = 111671406
$ws.Range("B9").Value = 78578
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 6458
$ws.Range("F9").Value = "Lunglav"
$ws.Range("G9").Value = "Lobaria pulmonaria"
$ws.Range("H9").Value = "(L.) Hoffm."
$ws.Range("Q9").Value = 557823
$ws.Range("R9").Value = 7068159
$ws.Range("Z9").ClearContents()
$ws.Range("AB9").ClearContents()

# Row 10
$ws.Range("A10").Value = 111670588
$ws.Range("B10").Value = 96348
$ws.Range("D10").Value = "VU"
$ws.Range("E10").Value = 220787
$ws.Range("F10").Value = "Knärot"
$ws.Range("G10").Value = "Goodyera repens"
$ws.Range("H10").Value = "(L.) R. Br."
$ws.Range("Q10").Value = 558040
$ws.Range("R10").Value = 7067902
$ws.Range("Z10").ClearContents()
$ws.Range("AB10").ClearContents()

# Row 11
$ws.Range("A11").Value = 111671345
$ws.Range("B11").Value = 96348
$ws.Range("D11").Value = "VU"
$ws.Range("E11").Value = 220787
$ws.Range("F11").Value = "Knärot"
$ws.Range("G11").Value = "Goodyera repens"
$ws.Range("H11").Value = "(L.) R. Br."
$ws.Range("Q11").Value = 557813
$ws.Range("R11").Value = 7068166
$ws.Range("Z11").ClearContents()
$ws.Range("AB11").ClearContents()

# Row 12
$ws.Range("A12").Value = 111671395
$ws.Range("B12").Value = 96348
$ws.Range("D12").Value = "VU"
$ws.Range("E12").Value = 220787
$ws.Range("F12").Value = "Knärot"
$ws.Range("G12").Value = "Goodyera repens"
$ws.Range("H12").Value = "(L.) R. Br."
$ws.Range("Q12").Value = 557763
$ws.Range("R12").Value = 7068265
$ws.Range("Z12").ClearContents()
$ws.Range("AB12").ClearContents()

# Row 13
$ws.Range("A13").Value = 111671384
$ws.Range("B13").Value = 96348
$ws.Range("D13").Value = "VU"
$ws.Range("E13").Value = 220787
$ws.Range("F13").Value = "Knärot"
$ws.Range("G13").Value = "Goodyera repens"
$ws.Range("H13").Value = "(L.) R. Br."
$ws.Range("Q13").Value = 557798
$ws.Range("R13").Value = 7068181
$ws.Range("Z13").ClearContents()
$ws.Range("AB13").ClearContents()

# Row 14
$ws.Range("A14").Value = 111670607
$ws.Range("B14").Value = 96368
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 221952
$ws.Range("F14").Value = "Spindelblomster"
$ws.Range("G14").Value = "Neottia cordata"
$ws.Range("H14").Value = "(L.) Rich."
$ws.Range("Q14").Value = 558032
$ws.Range("R14").Value = 7067908
$ws.Range("Z14").ClearContents()
$ws.Range("AB14").ClearContents()

# Row 15
$ws.Range("A15").Value = 111671364
$ws.Range("B15").Value = 96368
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 221952
$ws.Range("F15").Value = "Spindelblomster"
$ws.Range("G15").Value = "Neottia cordata"
$ws.Range("H15").Value = "(L.) Rich."
$ws.Range("Q15").Value = 557813
$ws.Range("R15").Value = 7068169
$ws.Range("Z15").ClearContents()
$ws.Range("AB15").ClearContents()

# Row 16
$ws.Range("A16").Value = 111670593
$ws.Range("B16").Value = 78578
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 6458
$ws.Range("F16").Value = "Lunglav"
$ws.Range("G16").Value = "Lobaria pulmonaria"
$ws.Range("H16").Value = "(L.) Hoffm."
$ws.Range("Q16").Value = 558041
$ws.Range("R16").Value = 7067901
$ws.Range("Z16").ClearContents()
$ws.Range("AB16").ClearContents()

# Row 17
$ws.Range("A17").Value = 111670599
$ws.Range("B17").Value = 96348
$ws.Range("D17").Value = "VU"
$ws.Range("E17").Value = 220787
$ws.Range("F17").Value = "Knärot"
$ws.Range("G17").Value = "Goodyera repens"
$ws.Range("H17").Value = "(L.) R. Br."
$ws.Range("Q17").Value = 558032
$ws.Range("R17").Value = 7067909
$ws.Range("Z17").ClearContents()
$ws.Range("AB17").ClearContents()

# Row 18
$ws.Range("A18").Value = 111670575
$ws.Range("B18").Value = 96346
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 620
$ws.Range("F18").Value = "Skogsfru"
$ws.Range("G18").Value = "Epipogium aphyllum"
$ws.Range("H18").Value = "Sw."
$ws.Range("Q18").Value = 558083
$ws.Range("R18").Value = 7067975
$ws.Range("Z18").ClearContents()
$ws.Range("AB18").ClearContents()
